# "10.02/2024 - return to host 33"
# Bump Date_of_sales (column I) from 45348 to 45349 for every data row (2-30),
# and reorder a handful of tokens in the Param / concatenated-Param columns
# (C and G) so that "б/к" / "Type" move to the front of their token lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Column I: Date_of_sales, rows 2..30, 45348 -> 45349 ----
for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 9).Value = 45349
}

# ---- Column C: reorder tokens ----
$ws.Cells.Item(3, 3).Value  = "б/к сер легк"
$ws.Cells.Item(4, 3).Value  = "б/к сер легк"
$ws.Cells.Item(10, 3).Value = "210B сер Type C H"
$ws.Cells.Item(11, 3).Value = "LS-2 сер груз Type"
$ws.Cells.Item(12, 3).Value = "202B сер Type C"
$ws.Cells.Item(13, 3).Value = "202B LS-2 сер Type C H"
$ws.Cells.Item(14, 3).Value = "б/к сер груз"
$ws.Cells.Item(15, 3).Value = "б/к сер легк"
$ws.Cells.Item(16, 3).Value = "б/к сер легк"

# ---- Column G: reorder tokens ----
$ws.Cells.Item(4, 7).Value  = "б/к, сер, легк"
$ws.Cells.Item(5, 7).Value  = "б/к, сер, легк"
$ws.Cells.Item(11, 7).Value = "210B, сер, Type, C, H"
$ws.Cells.Item(12, 7).Value = "210B, сер, Type, C, H"
$ws.Cells.Item(13, 7).Value = "LS-2, сер, груз, Type"
$ws.Cells.Item(14, 7).Value = "202B, сер, Type, C"
$ws.Cells.Item(15, 7).Value = "202B, LS-2, сер, Type, C, H"
$ws.Cells.Item(16, 7).Value = "202B, LS-2, сер, Type, C, H"
$ws.Cells.Item(17, 7).Value = "202B, LS-2, сер, Type, C, H"
$ws.Cells.Item(18, 7).Value = "б/к, сер, груз"
$ws.Cells.Item(19, 7).Value = "б/к, сер, груз"
$ws.Cells.Item(20, 7).Value = "б/к, сер, груз"
$ws.Cells.Item(21, 7).Value = "б/к, сер, груз"
$ws.Cells.Item(22, 7).Value = "б/к, сер, легк"
$ws.Cells.Item(23, 7).Value = "б/к, сер, легк"
